# Edit script: add "metadata" sheet, refresh time_taken timestamps on "data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. Refresh the time_taken values (column F) for rows 2..67 of "data" ---
$newTimes = @(
    "2021-10-05 14:34:11.659819",
    "2021-10-05 14:34:11.659827",
    "2021-10-05 14:34:11.659830",
    "2021-10-05 14:34:11.659833",
    "2021-10-05 14:34:11.659836",
    "2021-10-05 14:34:11.659839",
    "2021-10-05 14:34:11.659841",
    "2021-10-05 14:34:11.659844",
    "2021-10-05 14:34:11.659847",
    "2021-10-05 14:34:11.659850",
    "2021-10-05 14:34:11.659852",
    "2021-10-05 14:34:11.659855",
    "2021-10-05 14:34:11.659857",
    "2021-10-05 14:34:11.659860",
    "2021-10-05 14:34:11.659862",
    "2021-10-05 14:34:11.659865",
    "2021-10-05 14:34:11.659868",
    "2021-10-05 14:34:11.659870",
    "2021-10-05 14:34:11.659873",
    "2021-10-05 14:34:11.659876",
    "2021-10-05 14:34:11.659878",
    "2021-10-05 14:34:11.659881",
    "2021-10-05 14:34:11.659883",
    "2021-10-05 14:34:11.659886",
    "2021-10-05 14:34:11.659889",
    "2021-10-05 14:34:11.659892",
    "2021-10-05 14:34:11.659894",
    "2021-10-05 14:34:11.659897",
    "2021-10-05 14:34:11.659900",
    "2021-10-05 14:34:11.659902",
    "2021-10-05 14:34:11.659905",
    "2021-10-05 14:34:11.659907",
    "2021-10-05 14:34:11.659910",
    "2021-10-05 14:34:11.659913",
    "2021-10-05 14:34:11.659916",
    "2021-10-05 14:34:11.659918",
    "2021-10-05 14:34:11.659921",
    "2021-10-05 14:34:11.659923",
    "2021-10-05 14:34:11.659926",
    "2021-10-05 14:34:11.659929",
    "2021-10-05 14:34:11.659932",
    "2021-10-05 14:34:11.659934",
    "2021-10-05 14:34:11.659937",
    "2021-10-05 14:34:11.659940",
    "2021-10-05 14:34:11.659942",
    "2021-10-05 14:34:11.659945",
    "2021-10-05 14:34:11.659947",
    "2021-10-05 14:34:11.659950",
    "2021-10-05 14:34:11.659953",
    "2021-10-05 14:34:11.659955",
    "2021-10-05 14:34:11.659958",
    "2021-10-05 14:34:11.659961",
    "2021-10-05 14:34:11.659964",
    "2021-10-05 14:34:11.659966",
    "2021-10-05 14:34:11.659969",
    "2021-10-05 14:34:11.659971",
    "2021-10-05 14:34:11.659974",
    "2021-10-05 14:34:11.659977",
    "2021-10-05 14:34:11.659979",
    "2021-10-05 14:34:11.659982",
    "2021-10-05 14:34:11.659984",
    "2021-10-05 14:34:11.659987",
    "2021-10-05 14:34:11.659990",
    "2021-10-05 14:34:11.659993",
    "2021-10-05 14:34:11.659997",
    "2021-10-05 14:34:11.660000",
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add a new "metadata" worksheet placed right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row (A2:G2)
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Inflammatory bowel disease"
$meta.Range("C2").Value = 123

$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.61"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-09-20T09:17:23.391936Z"
$meta.Range("F2").Value = "2021-10-05 14:34:11.656531"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/123/?format=json"

# Apply the same bold/centered/bordered header style used on the "data" sheet
# to the new header row and the A2 index cell (reuse via copy/paste-format so
# no new style entries are introduced).
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Leave "data" as the active sheet, matching the original active tab.
$ws.Activate()
